$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 259.4
$ws.Range("I2").Value = 273.75
$ws.Range("J2").Value = 202
$ws.Range("K2").Value = 273.75
$ws.Range("L2").Value = 202
$ws.Range("M2").Value = -160.75
$ws.Range("N2").Value = -428
$ws.Range("H5").Value = 434.2
$ws.Range("J5").Value = 1050.5
$ws.Range("L5").Value = 1050.5
$ws.Range("N5").Value = -1280.5
$ws.Range("H11").Value = 166679.33
$ws.Range("I11").Value = 166679.33
$ws.Range("K11").Value = 166679.33
$ws.Range("M11").Value = -166539.33
$ws.Range("H19").Value = 851.5625
$ws.Range("I19").Value = 180.55556
$ws.Range("J19").Value = 1714.2858
$ws.Range("K19").Value = 180.55556
$ws.Range("L19").Value = 1714.2858
$ws.Range("M19").Value = -5.555560000000014
$ws.Range("N19").Value = -2064.2858
$ws.Range("H125").Value = 2113.739
$ws.Range("I125").Value = 2079.4285
$ws.Range("K125").Value = 18714.8565
$ws.Range("M125").Value = -16254.8565
$ws.Range("H129").Value = 2833.7874
$ws.Range("J129").Value = 830.71875
$ws.Range("L129").Value = 2492.15625
$ws.Range("N129").Value = -12492.15625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H32").Value = 6217.8706
$ws.Range("I32").Value = 5354.1357
$ws.Range("J32").Value = 23708.5
$ws.Range("K32").Value = 5354.1357
$ws.Range("L32").Value = 23708.5
$ws.Range("M32").Value = -5067.1357
$ws.Range("N32").Value = -24282.5
$ws.Range("H53").Value = 8266.666999999999
$ws.Range("I53").Value = 5000
$ws.Range("J53").Value = 9900
$ws.Range("K53").Value = 5000
$ws.Range("L53").Value = 9900
$ws.Range("M53").Value = -4318
$ws.Range("N53").Value = -11264

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H99").Value = 1295.0454
$ws.Range("I99").Value = 781.86957
$ws.Range("K99").Value = 781.86957
$ws.Range("M99").Value = 716.13043

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 129.35715
$ws.Range("I7").Value = 66.25
$ws.Range("J7").Value = 213.5
$ws.Range("K7").Value = 66.25
$ws.Range("L7").Value = 213.5
$ws.Range("M7").Value = 46.75
$ws.Range("N7").Value = -439.5
$ws.Range("H19").Value = 16595.834
$ws.Range("I19").Value = 37.5
$ws.Range("J19").Value = 24875
$ws.Range("K19").Value = 37.5
$ws.Range("L19").Value = 24875
$ws.Range("M19").Value = 132.5
$ws.Range("N19").Value = -25215
$ws.Range("H24").Value = 16595.834
$ws.Range("I24").Value = 37.5
$ws.Range("J24").Value = 24875
$ws.Range("K24").Value = 37.5
$ws.Range("L24").Value = 24875
$ws.Range("M24").Value = 132.5
$ws.Range("N24").Value = -25215
$ws.Range("H31").Value = 2946.7708
$ws.Range("I31").Value = 1721.7391
$ws.Range("K31").Value = 1721.7391
$ws.Range("M31").Value = -1426.7391
$ws.Range("H34").Value = 2946.7708
$ws.Range("I34").Value = 1721.7391
$ws.Range("K34").Value = 1721.7391
$ws.Range("M34").Value = -1519.7391

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 3210.8
$ws.Range("J117").Value = 3210.8
$ws.Range("L117").Value = 9632.400000000001
$ws.Range("N117").Value = -16516.4
$ws.Range("H121").Value = 6961.294
$ws.Range("I121").Value = 22844.5
$ws.Range("J121").Value = 4843.533
$ws.Range("K121").Value = 68533.5
$ws.Range("L121").Value = 14530.599
$ws.Range("M121").Value = -67223.5
$ws.Range("N121").Value = -17150.599
$ws.Range("H131").Value = 753.21
$ws.Range("J131").Value = 780.43616
$ws.Range("L131").Value = 2341.30848
$ws.Range("N131").Value = -12421.30848
$ws.Range("H132").Value = 1587
$ws.Range("I132").Value = 813.06665
$ws.Range("K132").Value = 7317.59985
$ws.Range("M132").Value = -4787.59985
$ws.Range("H133").Value = 3046.6667
$ws.Range("I133").Value = 2082.5
$ws.Range("J133").Value = 4975
$ws.Range("K133").Value = 6247.5
$ws.Range("L133").Value = 14925
$ws.Range("M133").Value = -1187.5
$ws.Range("N133").Value = -25045
$ws.Range("H136").Value = 2188.889
$ws.Range("I136").Value = 1783.3334
$ws.Range("K136").Value = 5350.0002
$ws.Range("M136").Value = -250.0002000000004
$ws.Range("H137").Value = 70784.56
$ws.Range("J137").Value = 28873.25
$ws.Range("L137").Value = 86619.75
$ws.Range("N137").Value = -96819.75

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 259.22223
$ws.Range("I2").Value = 122
$ws.Range("K2").Value = 122
$ws.Range("M2").Value = -9

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2033.5625
$ws.Range("I22").Value = 1534.2858
$ws.Range("K22").Value = 1534.2858
$ws.Range("M22").Value = -1239.2858
$ws.Range("H27").Value = 2033.5625
$ws.Range("I27").Value = 1534.2858
$ws.Range("K27").Value = 1534.2858
$ws.Range("M27").Value = -1427.2858
$ws.Range("H46").Value = 779065.4399999999
$ws.Range("I46").Value = 390.25
$ws.Range("J46").Value = 1125143.4
$ws.Range("K46").Value = 390.25
$ws.Range("L46").Value = 1125143.4
$ws.Range("M46").Value = -202.25
$ws.Range("N46").Value = -1125519.4
$ws.Range("H87").Value = 36389.668
$ws.Range("J87").Value = 36389.668
$ws.Range("L87").Value = 36389.668
$ws.Range("N87").Value = -38635.668
$ws.Range("H90").Value = 36389.668
$ws.Range("J90").Value = 36389.668
$ws.Range("L90").Value = 109169.004
$ws.Range("N90").Value = -120401.004
$ws.Range("H122").Value = 3047.4
$ws.Range("J122").Value = 3333
$ws.Range("L122").Value = 9999
$ws.Range("N122").Value = -14899
$ws.Range("H136").Value = 2137.5
$ws.Range("I136").Value = 2216.6667
$ws.Range("J136").Value = 1900
$ws.Range("K136").Value = 6650.000100000001
$ws.Range("L136").Value = 5700
$ws.Range("M136").Value = -4100.000100000001
$ws.Range("N136").Value = -10800
$ws.Range("H139").Value = 64618.332
$ws.Range("J139").Value = 64618.332
$ws.Range("L139").Value = 64618.332
$ws.Range("N139").Value = -74898.33199999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 12800
$ws.Range("J25").Value = 12800
$ws.Range("L25").Value = 12800
$ws.Range("N25").Value = -13386
$ws.Range("H122").Value = 2549.6667
$ws.Range("I122").Value = 1960.1666
$ws.Range("J122").Value = 2844.4167
$ws.Range("K122").Value = 5880.4998
$ws.Range("L122").Value = 8533.250100000001
$ws.Range("M122").Value = -3430.4998
$ws.Range("N122").Value = -13433.2501
$ws.Range("H138").Value = 69333.336
$ws.Range("J138").Value = 69333.336
$ws.Range("L138").Value = 69333.336
$ws.Range("N138").Value = -79613.336
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
